# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the second data row (the 4a22f876-... entry) on the zh-cn sheet,
# and for the second data row (the 4a22f876-... entry) on the de-de sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-20 22:42:21"
$wsZhCn.Range("H3").Value = "2016-03-20 22:42:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-20 22:42:24"
$wsDeDe.Range("H3").Value = "2016-03-20 22:42:53"
